# Adapt column header formatting to respective input file names.
# - Rename "<col>_old" headers to "<col>_FV2404"
# - Rename "<col>_new" headers to "<col>_FV2410"
# - Wrap the data range in an Excel Table ("Table1")
# - Freeze the header row (row 1) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404", "Segment ID_FV2404",
  "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404", "Bedingungsausdruck_FV2404", "Bedingung_FV2404",
  "diff",
  "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410", "Segment ID_FV2410",
  "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410", "Bedingungsausdruck_FV2410", "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into a proper Excel Table so the renamed headers
# become the table's column names.
$dataRange = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
